$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email text values in A4 and A7
$ws.Range("A4").Value = "rajumsuryawanshi19@gmail.com"
$ws.Range("A7").Value = "rajemsuryawanshi17@gmail.com"

# Move the active selection to A4 (was A8)
$ws.Range("A4").Select()
